# Edit LOQ4079.xlsx: populate previously-misaligned label/value pairs and
# insert a new row so every "Objetivos:/Objectives:/Docentes responsáveis:"
# trio has its own correctly-placed value row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 11 (everything from old row 11 downward shifts to row+1) ---
$ws.Rows.Item(11).Insert()

# Row 10 ("Objetivos:") now gets the real objectives text instead of the
# professor name that had been mistakenly placed there.
$ws.Cells.Item(10,2).Value = "Introduzir o aluno na engenharia das reações químicas, através dos conceitos fundamentais da cinética química aplicada a reatores químicos ideais."
$ws.Cells.Item(10,3).Value = "Introduzir o aluno na engenharia das reações químicas, através dos conceitos fundamentais da cinética química aplicada a reatores químicos ideais."

# New row 11 holds the "Objectives:" label (English) with the same row height
# as the "Objetivos:" row above it.
$ws.Cells.Item(11,1).Value = "Objectives:"
$ws.Rows.Item(11).RowHeight = 60

# Row 12 keeps the "Docentes responsáveis:" label (shifted down from 12->12,
# previously row 12); restore its default (non-custom) row height.
$ws.Cells.Item(12,1).Value = "Docentes responsáveis:"
$ws.Rows.Item(12).AutoFit()

# Row 13 (blank label row, shifted from old row 12 "Docentes responsáveis:")
# now carries the professor's name/id as its value.
$ws.Cells.Item(13,2).Value = "6310316 - Liana Alvares Rodrigues"
$ws.Cells.Item(13,3).Value = "6310316 - Liana Alvares Rodrigues"

# Row 14 "Programa resumido:" gets the real short-syllabus text instead of
# the placeholder "Semestral".
$ws.Cells.Item(14,2).Value = "1. Introdução a cinética. 2. Reações a volume constante. 3. Reações a volume variável. 4. Modelos ideais de reatores químicos isotérmicos. 5. Análise de dados cinéticos em reatores químicos isotérmicos"
$ws.Cells.Item(14,3).Value = "1. Introdução a cinética. 2. Reações a volume constante. 3. Reações a volume variável. 4. Modelos ideais de reatores químicos isotérmicos. 5. Análise de dados cinéticos em reatores químicos isotérmicos"

# Row 16 "Programa:" gets the full syllabus text instead of the duplicated
# activation date.
$ws.Cells.Item(16,2).Value = "1. INTRODUÇÃO A CINÉTICA`nTipos de Reações Químicas. Lei de velocidade e seus principais parâmetros. Influência da temperatura sobre a taxa da reação. Ativação das reações químicas Equação de Arrhenius. Energia de ativação. Conversão. Concentração e sua variação numa transformação química. (4 horas)`n2. REAÇÕES A VOLUME CONSTANTE`nReações irreversíveis de ordem um. Reações irreversíveis de ordem dois. Reações irreversíveis de ordem três. Reações irreversíveis de ordem qualquer. (8 horas)`n3. REAÇÕES A VOLUME VARIÁVEL`nConceitos. Fração de conversão volumétrica. Reações a volume variável de ordem um e dois. (2 horas)`n4. MODELOS IDEAIS DE REATORES QUÍMICOS ISOTÉRMICOS: `nEquações fundamentais de projeto de reatores. Reator tanque descontínuo (BSTR). Reator tanque de mistura contínuo (CSTR). Reator tubular de fluxo pistonado (PFR). Comparação de desempenho de reatores CSTR e PFR. Reatores CSTR em cascata. Associação mista de reatores em série: CSTR e PFR (8 horas)`n5. ANÁLISE DE DADOS CINÉTICOS EM REATORES QUÍMICOS ISOTÉRMICOS`nBalanço de massa e coleta de dados em reatores ideais isotérmicos: batelada (BSTR), reator tanque de mistura contínuo (CSTR) e Reator tubular (PFR) (8 horas)"
$ws.Cells.Item(16,3).Value = "1. INTRODUÇÃO A CINÉTICA`nTipos de Reações Químicas. Lei de velocidade e seus principais parâmetros. Influência da temperatura sobre a taxa da reação. Ativação das reações químicas Equação de Arrhenius. Energia de ativação. Conversão. Concentração e sua variação numa transformação química. (4 horas)`n2. REAÇÕES A VOLUME CONSTANTE`nReações irreversíveis de ordem um. Reações irreversíveis de ordem dois. Reações irreversíveis de ordem três. Reações irreversíveis de ordem qualquer. (8 horas)`n3. REAÇÕES A VOLUME VARIÁVEL`nConceitos. Fração de conversão volumétrica. Reações a volume variável de ordem um e dois. (2 horas)`n4. MODELOS IDEAIS DE REATORES QUÍMICOS ISOTÉRMICOS: `nEquações fundamentais de projeto de reatores. Reator tanque descontínuo (BSTR). Reator tanque de mistura contínuo (CSTR). Reator tubular de fluxo pistonado (PFR). Comparação de desempenho de reatores CSTR e PFR. Reatores CSTR em cascata. Associação mista de reatores em série: CSTR e PFR (8 horas)`n5. ANÁLISE DE DADOS CINÉTICOS EM REATORES QUÍMICOS ISOTÉRMICOS`nBalanço de massa e coleta de dados em reatores ideais isotérmicos: batelada (BSTR), reator tanque de mistura contínuo (CSTR) e Reator tubular (PFR) (8 horas)"

# Row 19 "Método:" gets the evaluation-method text instead of the professor
# name.
$ws.Cells.Item(19,2).Value = "Duas provas escritas (P1 e P2) e trabalhos relacionados à disciplina (TRAB)."
$ws.Cells.Item(19,3).Value = "Duas provas escritas (P1 e P2) e trabalhos relacionados à disciplina (TRAB)."

# Row 20 "Critério:" gets the grading-weights text.
$ws.Cells.Item(20,2).Value = "Média da Primeira Avaliação = (I)  Prova P1=30%; (II)  Prova P2=60% e (III)  Trabalhos =10%"
$ws.Cells.Item(20,3).Value = "Média da Primeira Avaliação = (I)  Prova P1=30%; (II)  Prova P2=60% e (III)  Trabalhos =10%"

# Row 21 "Norma de recuperação:" gets the make-up exam rule text.
$ws.Cells.Item(21,2).Value = "Será a média aritmética da nota do aluno na primeira avaliação e da nota do aluo numa prova escrita na recuperação."
$ws.Cells.Item(21,3).Value = "Será a média aritmética da nota do aluno na primeira avaliação e da nota do aluo numa prova escrita na recuperação."

# Row 22 "Bibliografia:" gets the actual bibliography instead of the
# make-up exam rule text.
$ws.Cells.Item(22,2).Value = "FOGLER, H. S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2002.`nLEVENSPIEL, O. Chemical Reaction Engineering. 3rd. ed. New York: John Wiley & Sons, 1998.`nHILL, C.G. An Introduction to chemical engineering kinetics and reactor design. New York: John Wiley&Sons, 1977.`nSMITH, J.M. Chemical engineering kinetics. 3rd. ed New York: McGraw-Hill,1981.`nDENBIGH, K. ; TURNER, R. Introduction to chemical Reaction Design. Cambridge: Cambridge University Press, 1970.`nFROMENT, G.F. ; BISCHOFF, K.B. Chemical reactor analysis and design. 2nd Ed. New York: John Wiley & Sons, 1990."
$ws.Cells.Item(22,3).Value = "FOGLER, H. S. Elementos de engenharia das reações químicas. 3.ed. Rio de Janeiro: LTC Editora, 2002.`nLEVENSPIEL, O. Chemical Reaction Engineering. 3rd. ed. New York: John Wiley & Sons, 1998.`nHILL, C.G. An Introduction to chemical engineering kinetics and reactor design. New York: John Wiley&Sons, 1977.`nSMITH, J.M. Chemical engineering kinetics. 3rd. ed New York: McGraw-Hill,1981.`nDENBIGH, K. ; TURNER, R. Introduction to chemical Reaction Design. Cambridge: Cambridge University Press, 1970.`nFROMENT, G.F. ; BISCHOFF, K.B. Chemical reactor analysis and design. 2nd Ed. New York: John Wiley & Sons, 1990."
